$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.02635514705930575
$ws.Range("D2").Value = 0.2674911136138434
$ws.Range("E2").Value = 0.1853224059509699
$ws.Range("F2").Value = 0.9284568350695821
$ws.Range("G2").Value = 0.4121690984244708
$ws.Range("H2").Value = 0.5826480381023345
$ws.Range("I2").Value = 0.4268648828551349
$ws.Range("J2").Value = 0.1743904219952839
$ws.Range("K2").Value = 1.828067509145853
$ws.Range("O2").Value = 1.93774192686999
$ws.Range("C3").Value = 0.02311475358142445
$ws.Range("D3").Value = 0.2590922776626599
$ws.Range("E3").Value = 0.1804418432915682
$ws.Range("F3").Value = 0.9362983061392853
$ws.Range("G3").Value = 0.4191390201237724
$ws.Range("H3").Value = 0.5906760344763242
$ws.Range("I3").Value = 0.4319089911625795
$ws.Range("J3").Value = 0.170554697630152
$ws.Range("K3").Value = 1.60584553751039
$ws.Range("O3").Value = 1.969035205413377
$ws.Range("C4").Value = 0.02111596673712768
$ws.Range("D4").Value = 0.2540132657260159
$ws.Range("E4").Value = 0.1775288732370797
$ws.Range("F4").Value = 0.9418380448627346
$ws.Range("G4").Value = 0.4238785257555975
$ws.Range("H4").Value = 0.5959741750421301
$ws.Range("I4").Value = 0.4353982354070176
$ws.Range("J4").Value = 0.1683040952776551
$ws.Range("K4").Value = 1.468850518247905
$ws.Range("O4").Value = 1.989987087726519
$ws.Range("C5").Value = 0.02029918221253979
$ws.Range("D5").Value = 0.2519632894279908
$ws.Range("E5").Value = 0.1763629303569552
$ws.Range("F5").Value = 0.944277579976081
$ws.Range("G5").Value = 0.425925187403891
$ws.Range("H5").Value = 0.5982258901971633
$ws.Range("I5").Value = 0.4369185196431538
$ws.Range("J5").Value = 0.1674132496256888
$ws.Range("K5").Value = 1.412889922205409
$ws.Range("O5").Value = 1.998961104000045
$ws.Range("C6").Value = 0.02016342040002428
$ws.Range("D6").Value = 0.2516240910952945
$ws.Range("E6").Value = 0.1761706037169688
$ws.Range("F6").Value = 0.9446936501080003
$ws.Range("G6").Value = 0.4262719850761343
$ws.Range("H6").Value = 0.598605380843459
$ws.Range("I6").Value = 0.4371768977107031
$ws.Range("J6").Value = 0.1672669135377163
$ws.Range("K6").Value = 1.403589739590075
$ws.Range("O6").Value = 2.000477533725061
$ws.Range("C7").Value = 0.02110496039835397
$ws.Range("D7").Value = 0.2539855387830983
$ws.Range("E7").Value = 0.1775130633382602
$ws.Range("F7").Value = 0.9418702085487638
$ws.Range("G7").Value = 0.4239056615164074
$ws.Range("H7").Value = 0.5960041673207641
$ws.Range("I7").Value = 0.4354183404097647
$ws.Range("J7").Value = 0.1682919745568014
$ws.Range("K7").Value = 1.468096351034205
$ws.Range("O7").Value = 1.990106350619172
$ws.Range("C8").Value = 0.02523978827341011
$ws.Range("D8").Value = 0.2645791251590168
$ws.Range("E8").Value = 0.1836222499935829
$ws.Range("F8").Value = 0.9310099277283257
$ws.Range("G8").Value = 0.4144766365824637
$ws.Range("H8").Value = 0.5853394702521513
$ws.Range("I8").Value = 0.4285225714354368
$ws.Range("J8").Value = 0.1730461658246227
$ws.Range("K8").Value = 1.751561957448303
$ws.Range("O8").Value = 1.948170594909286
$ws.Range("C9").Value = 0.03327385993537746
$ws.Range("D9").Value = 0.2859649617225983
$ws.Range("E9").Value = 0.1962645314601659
$ws.Range("F9").Value = 0.9154785840806738
$ws.Range("G9").Value = 0.3996529154313109
$ws.Range("H9").Value = 0.5673570088165363
$ws.Range("I9").Value = 0.4181213868384326
$ws.Range("J9").Value = 0.1831992333157473
$ws.Range("K9").Value = 2.302915645758787
$ws.Range("O9").Value = 1.879766889370799
$ws.Range("C10").Value = 0.03912969255108578
$ws.Range("D10").Value = 0.3020434721404115
$ws.Range("E10").Value = 0.2059549304365262
$ws.Range("F10").Value = 0.9076009276289838
$ws.Range("G10").Value = 0.3910206343484717
$ws.Range("H10").Value = 0.5559370687278644
$ws.Range("I10").Value = 0.4123963961665389
$ws.Range("J10").Value = 0.1911666062996886
$ws.Range("K10").Value = 2.705065989585762
$ws.Range("O10").Value = 1.838003714880259
$ws.Range("C11").Value = 0.04178322148938207
$ws.Range("D11").Value = 0.3094362617532056
$ws.Range("E11").Value = 0.2104504040774486
$ws.Range("F11").Value = 0.9047883527573219
$ws.Range("G11").Value = 0.3875890519193135
$ws.Range("H11").Value = 0.5511319700205064
$ws.Range("I11").Value = 0.4102111301669211
$ws.Range("J11").Value = 0.1949019220324431
$ws.Range("K11").Value = 2.887345073472886
$ws.Range("O11").Value = 1.820861670184556
$ws.Range("C12").Value = 0.04278652420896378
$ws.Range("D12").Value = 0.3122468745402784
$ws.Range("E12").Value = 0.2121652238803406
$ws.Range("F12").Value = 0.9038344834192316
$ws.Range("G12").Value = 0.3863612315382241
$ws.Range("H12").Value = 0.5493685664740866
$ws.Range("I12").Value = 0.4094441256193662
$ws.Range("J12").Value = 0.1963323532053778
$ws.Range("K12").Value = 2.956270744208098
$ws.Range("O12").Value = 1.814638424017105
$ws.Range("C13").Value = 0.04257051386316846
$ws.Range("D13").Value = 0.311641066879389
$ws.Range("E13").Value = 0.2117953526970382
$ws.Range("F13").Value = 0.9040349658738194
$ws.Range("G13").Value = 0.3866224714487601
$ws.Range("H13").Value = 0.5497458463568208
$ws.Range("I13").Value = 0.4096066187224778
$ws.Range("J13").Value = 0.1960235749881321
$ws.Range("K13").Value = 2.941430860824482
$ws.Range("O13").Value = 1.815966770501106
$ws.Range("C14").Value = 0.04186579476382235
$ws.Range("D14").Value = 0.3096672706377035
$ws.Range("E14").Value = 0.2105912335686426
$ws.Range("F14").Value = 0.9047076470725983
$ws.Range("G14").Value = 0.3874866000961887
$ws.Range("H14").Value = 0.5509857670861535
$ws.Range("I14").Value = 0.4101468141518865
$ws.Range("J14").Value = 0.1950192848511421
$ws.Range("K14").Value = 2.893017648252396
$ws.Range("O14").Value = 1.820344300803498
$ws.Range("C15").Value = 0.04143393317407629
$ws.Range("D15").Value = 0.3084597068285007
$ws.Range("E15").Value = 0.2098552993154996
$ws.Range("F15").Value = 0.9051341736728062
$ws.Range("G15").Value = 0.3880252467370013
$ws.Range("H15").Value = 0.5517525745057696
$ws.Range("I15").Value = 0.410485586738524
$ws.Range("J15").Value = 0.1944062050137489
$ws.Range("K15").Value = 2.86335007515811
$ws.Range("O15").Value = 1.823060609719562
$ws.Range("C16").Value = 0.0389560665064721
$ws.Range("D16").Value = 0.3015619022394276
$ws.Range("E16").Value = 0.2056628908979334
$ws.Range("F16").Value = 0.9078002782371186
$ws.Range("G16").Value = 0.3912548990310256
$ws.Range("H16").Value = 0.5562589488464482
$ws.Range("I16").Value = 0.4125476622370705
$ws.Range("J16").Value = 0.1909247264431144
$ws.Range("K16").Value = 2.693139942567882
$ws.Range("O16").Value = 1.839161438401888
$ws.Range("C17").Value = 0.0374332989572963
$ws.Range("D17").Value = 0.2973503255611831
$ws.Range("E17").Value = 0.2031132866159027
$ws.Range("F17").Value = 0.9096335507090814
$ws.Range("G17").Value = 0.3933633429418322
$ws.Range("H17").Value = 0.5591234072278795
$ws.Range("I17").Value = 0.4139201871672853
$ws.Range("J17").Value = 0.1888173596846201
$ws.Range("K17").Value = 2.588549067818633
$ws.Range("O17").Value = 1.849515021074282
$ws.Range("C18").Value = 0.03655647490222691
$ws.Range("D18").Value = 0.294935346326497
$ws.Range("E18").Value = 0.2016550392133141
$ws.Range("F18").Value = 0.9107605526423015
$ws.Range("G18").Value = 0.3946226575208414
$ws.Range("H18").Value = 0.5608076591536459
$ws.Range("I18").Value = 0.4147490564950154
$ws.Range("J18").Value = 0.1876156982855264
$ws.Range("K18").Value = 2.52832919429369
$ws.Range("O18").Value = 1.855644776978323
$ws.Range("C19").Value = 0.03625943245916119
$ws.Range("D19").Value = 0.2941189536696243
$ws.Range("E19").Value = 0.2011627151062569
$ws.Range("F19").Value = 0.9111545871980127
$ws.Range("G19").Value = 0.3950570293104789
$ws.Range("H19").Value = 0.5613842160870846
$ws.Range("I19").Value = 0.415036461850061
$ws.Range("J19").Value = 0.1872106295490283
$ws.Range("K19").Value = 2.507929271332841
$ws.Range("O19").Value = 1.857750168012871
$ws.Range("C20").Value = 0.03759550073576179
$ws.Range("D20").Value = 0.297797889651946
$ws.Range("E20").Value = 0.2033838463436481
$ws.Range("F20").Value = 0.9094308845189474
$ws.Range("G20").Value = 0.393134070166731
$ws.Range("H20").Value = 0.5588146824861226
$ws.Range("I20").Value = 0.4137699969610438
$ws.Range("J20").Value = 0.1890406119620565
$ws.Range("K20").Value = 2.599689390494575
$ws.Range("O20").Value = 1.848394780656463
$ws.Range("C21").Value = 0.04207282977532145
$ws.Range("D21").Value = 0.3102467219571281
$ws.Range("E21").Value = 0.2109445743318474
$ws.Range("F21").Value = 0.9045070440373522
$ws.Range("G21").Value = 0.3872308367617094
$ws.Range("H21").Value = 0.5506200465857702
$ws.Range("I21").Value = 0.4099865015161583
$ws.Range("J21").Value = 0.1953138363941775
$ws.Range("K21").Value = 2.907240516489026
$ws.Range("O21").Value = 1.819051229446416
$ws.Range("C22").Value = 0.04499007695868329
$ws.Range("D22").Value = 0.3184475095926018
$ws.Range("E22").Value = 0.2159586680042338
$ws.Range("F22").Value = 0.9019372876199867
$ws.Range("G22").Value = 0.3837905319907264
$ws.Range("H22").Value = 0.5455919103711082
$ws.Range("I22").Value = 0.4078665708222218
$ws.Range("J22").Value = 0.1995067084435789
$ws.Range("K22").Value = 3.107661209186745
$ws.Range("O22").Value = 1.801436593587695
$ws.Range("C23").Value = 0.0434339213552164
$ws.Range("D23").Value = 0.3140647282847624
$ws.Range("E23").Value = 0.2132759198970007
$ws.Range("F23").Value = 0.9032493937534269
$ws.Range("G23").Value = 0.3855883191837606
$ws.Range("H23").Value = 0.54824551338276
$ws.Range("I23").Value = 0.4089656548414737
$ws.Range("J23").Value = 0.1972603887279689
$ws.Range("K23").Value = 3.000747582457223
$ws.Range("O23").Value = 1.810694456510006
$ws.Range("C24").Value = 0.03752217355946641
$ws.Range("D24").Value = 0.2975955262523371
$ws.Range("E24").Value = 0.2032615027566536
$ws.Range("F24").Value = 0.9095222824460976
$ws.Range("G24").Value = 0.3932375775756469
$ws.Range("H24").Value = 0.5589541403364393
$ws.Range("I24").Value = 0.4138377740395676
$ws.Range("J24").Value = 0.1889396487962642
$ws.Range("K24").Value = 2.594653128088339
$ws.Range("O24").Value = 1.84890068901953
$ws.Range("C25").Value = 0.03110854235949034
$ws.Range("D25").Value = 0.2801147162765716
$ws.Range("E25").Value = 0.1927737310386419
$ws.Range("F25").Value = 0.9190610700407973
$ws.Range("G25").Value = 0.4032681602355623
$ws.Range("H25").Value = 0.5719074448379331
$ws.Range("I25").Value = 0.42059957627351
$ws.Range("J25").Value = 0.1803634979509923
$ws.Range("K25").Value = 2.154261976678924
$ws.Range("O25").Value = 1.896784796157604
